# Updates cryptos price/volume table cells (columns D and E) to match
# the latest scrape, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.794.02"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.351.03"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'544.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'136.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").Value = "2.349.56"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "'5.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'24.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "2.774.38"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "60.728.82"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "2.352.13"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'10.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'320.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'63.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  -8.26%  "
$ws.Range("D26").Value = "'8.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.81%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "'8.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'1.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("D30").Value = "'495.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.61%  "
$ws.Range("D31").Value = "0.0₃0867"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'4.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'18.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("E39").Value = "  +6.01%  "
$ws.Range("D40").Value = "'5.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("D41").Value = "'144.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.83%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'142.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'2.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.55%  "
$ws.Range("D46").Value = "'0.0517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'19.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.64%  "
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "'0.0903"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Value = "'11.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
